# Swap the full content of row 10 and row 11 (columns A..AY) in the
# active worksheet. Both rows describe a single observation record;
# this update reorders them (row 10 <-> row 11) while every other row
# stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 10
$row2 = 11
$firstCol = 1    # A
$lastCol  = 51   # AY

# Capture both rows' current values (and whether each cell actually
# holds content) before writing anything, since row1/row2 get
# overwritten from each other below.
$vals1 = @()
$vals2 = @()
$present1 = @()
$present2 = @()

for ($c = $firstCol; $c -le $lastCol; $c++) {
    $v1 = $ws.Cells.Item($row1, $c).Value2
    $v2 = $ws.Cells.Item($row2, $c).Value2
    $vals1 += , $v1
    $vals2 += , $v2
    $present1 += , ($v1 -ne $null)
    $present2 += , ($v2 -ne $null)
}

function Needs-TextForce([string]$value) {
    # Patterns Excel's cell-value setter auto-coerces away from plain
    # text: a bare number, or an ISO yyyy-mm-dd date. An empty string
    # also needs forcing, since a plain "" assignment clears the cell
    # instead of leaving a present-but-empty text cell.
    if ($value -eq "") { return $true }
    if ($value -match '^-?\d+(\.\d+)?$') { return $true }
    if ($value -match '^\d{4}-\d{2}-\d{2}$') { return $true }
    return $false
}

function Values-Equal($a, $b) {
    if ($a -eq $null -or $b -eq $null) { return $a -eq $b }
    if (($a -is [string]) -ne ($b -is [string])) { return $false }
    return $a -eq $b
}

function Set-SwappedCell($ws, $row, $col, $present, $value, $alreadyPresent, $alreadyValue) {
    $cell = $ws.Cells.Item($row, $col)

    # No-op when the destination already holds exactly this content
    # (e.g. both rows had the same empty cell) - avoids touching the
    # cell's stored representation/style for a change that isn't one.
    if (($present -eq $alreadyPresent) -and (Values-Equal $value $alreadyValue)) {
        return
    }

    if (-not $present) {
        $cell.ClearContents()
        return
    }

    if ($value -is [string]) {
        if (Needs-TextForce $value) {
            # Leading apostrophe forces text interpretation, so Excel
            # doesn't silently turn it into a number/date, and an
            # empty string still leaves a real (blank) text cell
            # behind instead of clearing it. Re-normalising the style
            # afterwards drops the "quote prefix" marker Excel tacks
            # on, so the cell ends up plain text with default
            # formatting (matching how the source file stored it).
            $cell.Value = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    } else {
        $cell.Value = $value
    }
}

for ($i = 0; $i -lt ($lastCol - $firstCol + 1); $i++) {
    $c = $firstCol + $i
    Set-SwappedCell $ws $row1 $c $present2[$i] $vals2[$i] $present1[$i] $vals1[$i]
    Set-SwappedCell $ws $row2 $c $present1[$i] $vals1[$i] $present2[$i] $vals2[$i]
}
